$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Symbol"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Current Price"
$ws.Range("D1").Value = "Sector"
$ws.Range("E1").Value = "Description"
$ws.Range("F1").Value = "Risk"
$ws.Range("G1").Value = "Circuit Limit"
$ws.Range("H1").Value = "Volatility Factor"

# ---- Row 2: INFY ----
$ws.Range("A2").Value = "INFY"
$ws.Range("B2").Value = "Infosys Ltd."
$ws.Range("C2").Value = 1450
$ws.Range("D2").Value = "IT Services"
$ws.Range("E2").Value = "Leading IT consulting and outsourcing firm"
$ws.Range("F2").Value = "Medium"
$ws.Range("G2").Value = 100
$ws.Range("H2").Value = 60

# ---- Row 3: RELI ----
$ws.Range("A3").Value = "RELI"
$ws.Range("B3").Value = "Reliance Ind."
$ws.Range("C3").Value = 2750
$ws.Range("D3").Value = "Conglomerate"
$ws.Range("E3").Value = "Diversified biz in energy, retail, telecom"
$ws.Range("F3").Value = "Medium"
$ws.Range("G3").Value = 150
$ws.Range("H3").Value = 70

# ---- Row 4: TCS ----
$ws.Range("A4").Value = "TCS"
$ws.Range("B4").Value = "Tata Consultancy"
$ws.Range("C4").Value = 3600
$ws.Range("D4").Value = "IT Services"
$ws.Range("E4").Value = "Top global IT services and consulting firm"
$ws.Range("F4").Value = "Low"
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = 40

# ---- Row 5: ICIC ----
$ws.Range("A5").Value = "ICIC"
$ws.Range("B5").Value = "ICICI Bank"
$ws.Range("C5").Value = 1120
$ws.Range("D5").Value = "Banking"
$ws.Range("E5").Value = "Major private sector bank in India"
$ws.Range("F5").Value = "Low"
$ws.Range("G5").Value = 75
$ws.Range("H5").Value = 55

# ---- Row 6: DMART ----
$ws.Range("A6").Value = "DMART"
$ws.Range("B6").Value = "Avenue Supermarts"
$ws.Range("C6").Value = 4100
$ws.Range("D6").Value = "Retail"
$ws.Range("E6").Value = "Operates D-Mart chain across India"
$ws.Range("F6").Value = "High"
$ws.Range("G6").Value = 200
$ws.Range("H6").Value = 85

# ---- Header formatting: bold, centered, wrapped ----
$header = $ws.Range("A1:H1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108
$header.WrapText = $true

# ---- Body formatting: vertically centered, wrapped ----
$body = $ws.Range("A2:H6")
$body.VerticalAlignment = -4108
$body.WrapText = $true

# ---- Row heights ----
$ws.Rows.Item(1).RowHeight = 28.8
$ws.Rows.Item(2).RowHeight = 86.4
$ws.Rows.Item(3).RowHeight = 72
$ws.Rows.Item(4).RowHeight = 86.4
$ws.Rows.Item(5).RowHeight = 72
$ws.Rows.Item(6).RowHeight = 72

# ---- Column widths ----
$ws.Columns.Item(1).ColumnWidth = 7.33203125
$ws.Columns.Item(2).ColumnWidth = 21.109375
$ws.Columns.Item(3).ColumnWidth = 7.33203125
$ws.Columns.Item(4).ColumnWidth = 15.77734375
$ws.Columns.Item(5).ColumnWidth = 23.88671875

# ---- Selection ----
$ws.Range("J3").Select()
